$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "73.069.49"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3
$ws.Range("D3").Value = "3.980.94"
$ws.Range("E3").Value = "  -1.62%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.91%  "

# Row 7
$ws.Range("E7").Value = "  +0.29%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.802"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.74%  "

# Row 10
$ws.Range("E10").Value = "  +9.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.93%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.80%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.43%  "

# Row 14
$ws.Range("D14").Value = "4.618.28"
$ws.Range("E14").Value = "  -1.62%  "

# Row 15
$ws.Range("D15").Value = "3.972.33"
$ws.Range("E15").Value = "  -1.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.13%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "

# Row 18
$ws.Range("E18").Value = "  +0.88%  "

# Row 19
$ws.Range("D19").Value = "73.031.74"
$ws.Range("E19").Value = "  -0.14%  "

# Row 20
$ws.Range("E20").Value = "  -0.73%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "457.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.42%  "

# Row 22
$ws.Range("E22").Value = "  +6.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "97.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.15%  "

# Row 25
$ws.Range("E25").Value = "  -1.52%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.09%  "

# Row 29
$ws.Range("E29").Value = "  -1.40%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.47%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "14.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.82%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.07%  "

# Row 35
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0000103"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +16.71%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.37%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "634.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.97%  "

# Row 38
$ws.Range("E38").Value = "  -3.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.39%  "

# Row 40
$ws.Range("E40").Value = "  -1.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +47.48%  "

# Row 44
$ws.Range("E44").Value = "  -1.59%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.55%  "

# Row 46
$ws.Range("E46").Value = "  -0.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.19%  "

# Row 48
$ws.Range("E48").Value = "  +9.48%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.74%  "

# Row 50
$ws.Range("E50").Value = "  +1.78%  "

# Row 51
$ws.Range("E51").Value = "  -1.60%  "
